$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8 ("Upload"), pushing it (and
# everything below) down by one. The new row becomes row 8, "Upload" becomes
# row 9, the header row becomes row 10, and the TALENT1..8 data rows become
# rows 11..18.
$ws.Rows("8:8").Insert()

# Match the boxed / shaded style used by the surrounding boolean-flag rows
# (rows 3-7 and, after the shift, row 9) by giving the new row the same thin
# border used by that style.
$ws.Range("A8:G8").Borders.LineStyle = 1

# Populate the new row: a label in column A and FALSE flags in B:G, just like
# every other row in this block (Public/Private/Save/Cache/Ref/Upload).
$ws.Range("A8").Value = "Force"
$ws.Range("B8:G8").Value = $false

# Re-point the frozen pane / selection to follow the inserted row: the split
# now sits below row 10 (was row 9) and the active cell moves to A9 (the
# "Upload" row, now one lower than before).
$win = $excel.ActiveWindow
[void]($win.FreezePanes = $false)
[void]$ws.Range("A11").Select()
[void]($win.FreezePanes = $true)
[void]$ws.Range("A9").Select()
